# Adds a new "13.0.0" release column (L) to the project version matrix,
# mirroring the existing per-project version columns, then moves the
# active selection onto the newly added data (L5) as in the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tardigrade")

# Header for the new release column.
$ws.Range("L1").Value = "13.0.0"
$ws.Range("L1").Font.Bold = $true

# Per-project version values for the new release (blanks left empty,
# matching rows that have no release under this version).
$ws.Range("L2").Value  = "11.0.0"
$ws.Range("L4").Value  = "4.0.0"
$ws.Range("L7").Value  = "11.1.0"
$ws.Range("L8").Value  = "10.0.0"

# Rows with no release under this version still get a (blank) cell in the
# new column, matching the styled-but-empty cells already used in the
# equivalent spots of column K.
$ws.Range("L3").Value  = ""
$ws.Range("L3").Font.Bold = $false
$ws.Range("L6").Value  = ""
$ws.Range("L6").Font.Bold = $false
$ws.Range("L9").Value  = ""
$ws.Range("L9").Font.Bold = $false
$ws.Range("L10").Value = ""
$ws.Range("L10").Font.Bold = $false
$ws.Range("L11").Value = ""
$ws.Range("L11").Font.Bold = $false

# Move the selection to L5, matching the post-edit active cell.
$ws.Range("L5").Select()
